$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 - new timestamp, forecast columns only (M/N/O left blank)
$ws.Range("A19").Value = "2017.05.30 02.59.50"
$ws.Range("B19").Value = 15.14999961853027
$ws.Range("C19").Value = 100
$ws.Range("D19").Value = 1013
$ws.Range("E19").Value = 1.5
$ws.Range("F19").Value = "2017-05-30T09:00:00"
$ws.Range("G19").Value = "2017-05-30T12:00:00"
$ws.Range("H19").Value = 15.07999992370605
$ws.Range("I19").Value = 994.3699951171875
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 0.05999999865889549
$ws.Range("L19").Value = 1.610000014305115

# Row 20 - new timestamp, includes current-sensor columns M/N/O
$ws.Range("A20").Value = "2017.05.30 03.00.27"
$ws.Range("B20").Value = 15.14999961853027
$ws.Range("C20").Value = 100
$ws.Range("D20").Value = 1013
$ws.Range("E20").Value = 1.5
$ws.Range("F20").Value = "2017-05-30T09:00:00"
$ws.Range("G20").Value = "2017-05-30T12:00:00"
$ws.Range("H20").Value = 15.07999992370605
$ws.Range("I20").Value = 994.3699951171875
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 0.05999999865889549
$ws.Range("L20").Value = 1.610000014305115
$ws.Range("M20").Value = 16.98
$ws.Range("N20").Value = 2.33
$ws.Range("O20").Value = 84.2

# Row 21 - new timestamp, includes current-sensor columns M/N/O
$ws.Range("A21").Value = "2017.05.30 03.00.48"
$ws.Range("B21").Value = 15.14999961853027
$ws.Range("C21").Value = 100
$ws.Range("D21").Value = 1013
$ws.Range("E21").Value = 1.5
$ws.Range("F21").Value = "2017-05-30T09:00:00"
$ws.Range("G21").Value = "2017-05-30T12:00:00"
$ws.Range("H21").Value = 15.07999992370605
$ws.Range("I21").Value = 994.3699951171875
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 0.05999999865889549
$ws.Range("L21").Value = 1.610000014305115
$ws.Range("M21").Value = 17
$ws.Range("N21").Value = 2.34
$ws.Range("O21").Value = 84
